# New crime data collected — weekly CompStat figures roll forward one week.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume/issue number and the "week covering" date range ---
# These are rich-text shared strings; re-assign the full display text.
$ws.Range("A8").Value = "Volume 30   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/23/2023  Through  10/29/2023"

# --- Bulk numeric updates: same cell style/number-format, value changes only ---
$values = @{
    "G15" = 2
    "H15" = 0
    "J15" = 18
    "K15" = 22.222222222222
    "C16" = 7
    "F16" = 24
    "G16" = 13
    "H16" = 84.615384615384
    "I16" = 219
    "J16" = 161
    "K16" = 36.024844720496
    "L16" = 123.469387755102
    "M16" = 37.735849056603
    "N16" = -74.082840236686
    "C17" = 4
    "D17" = 2
    "E17" = 100
    "F17" = 19
    "G17" = 13
    "H17" = 46.153846153846
    "I17" = 205
    "J17" = 186
    "K17" = 10.215053763440
    "L17" = 17.142857142857
    "M17" = 99.029126213592
    "N17" = -22.053231939163
    "C18" = 5
    "D18" = 2
    "E18" = 150
    "F18" = 16
    "G18" = 12
    "H18" = 33.333333333333
    "I18" = 176
    "J18" = 139
    "K18" = 26.618705035971
    "L18" = 17.333333333333
    "M18" = -13.725490196078
    "N18" = -86.217697729052
    "C19" = 7
    "D19" = 12
    "E19" = -41.666666666666
    "F19" = 53
    "G19" = 53
    "H19" = 0
    "I19" = 621
    "J19" = 587
    "K19" = 5.792163543441
    "L19" = 68.292682926829
    "M19" = 61.298701298701
    "N19" = -18.823529411764
    "C20" = 5
    "D20" = 4
    "E20" = 25
    "F20" = 38
    "G20" = 23
    "H20" = 65.217391304347
    "I20" = 261
    "J20" = 196
    "K20" = 33.163265306122
    "L20" = 68.387096774193
    "M20" = 52.631578947368
    "N20" = -84.965437788018
    "C21" = 28
    "D21" = 24
    "E21" = 16.666666666666
    "F21" = 152
    "G21" = 116
    "H21" = 31.034482758620
    "I21" = 1504
    "J21" = 1288
    "K21" = 16.770186335403
    "L21" = 56.666666666666
    "M21" = 45.033751205400
    "N21" = -69.318645450836
    "C22" = 1
    "D22" = 2
    "E22" = -50
    "I22" = 71
    "J22" = 68
    "K22" = 4.411764705882
    "L22" = 108.823529411765
    "M22" = 61.363636363636
    "C24" = 39
    "D24" = 39
    "E24" = 0
    "F24" = 184
    "G24" = 144
    "H24" = 27.777777777777
    "I24" = 1711
    "J24" = 1249
    "K24" = 36.989591673338
    "L24" = 61.262959472196
    "M24" = 117.131979695431
    "C25" = 13
    "D25" = 13
    "E25" = 0
    "F25" = 47
    "G25" = 63
    "H25" = -25.396825396825
    "I25" = 436
    "J25" = 469
    "K25" = -7.036247334754
    "L25" = 10.101010101010
    "M25" = 2.830188679245
    "E26" = -100
    "G26" = 3
    "H26" = 33.333333333333
    "J26" = 24
    "K26" = 20.833333333333
    "D27" = 3
    "E27" = -66.666666666666
    "G27" = 9
    "H27" = -33.333333333333
    "I27" = 70
    "J27" = 76
    "K27" = -7.894736842105
    "L27" = 34.615384615384
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# --- Cells flipping from the "no data" text placeholder to real numbers ---
# (previously showed the shared "0"/"***.*" labels; now have counts)
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'

$ws.Range("D16").Value = 3
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("E16").Value = 133.333333333333
$ws.Range("E16").NumberFormat = '#,##0.0;"-"#,##0.0'

# --- Cells flipping from real numbers back to the "no data" text placeholder ---
# Copy formatting+value from the existing "0" / "***.*" template cells on row 23
# so the shared-string reference and style index line up exactly.
$ws.Range("C23").Copy($ws.Range("C26"))
$ws.Range("C23").Copy($ws.Range("C28"))
$ws.Range("C23").Copy($ws.Range("G28"))
$ws.Range("E23").Copy($ws.Range("H28"))
$ws.Range("C23").Copy($ws.Range("C29"))
$ws.Range("C23").Copy($ws.Range("G29"))
$ws.Range("E23").Copy($ws.Range("H29"))
